$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new daily rows (03-08-2021 and 04-08-2021) after the last existing row (147)
# Force column A to be stored as plain text (matching existing "dd-mm-yyyy" string labels)
# instead of letting Excel auto-convert it to a date serial number.
$ws.Range("A148:A149").NumberFormat = "@"

$ws.Range("A148").Value = "03-08-2021"
$ws.Range("B148").Value = -9764
$ws.Range("C148").Value = 3362
$ws.Range("D148").Value = 987
$ws.Range("E148").Value = 827
$ws.Range("F148").Value = 1549

$ws.Range("A149").Value = "04-08-2021"
$ws.Range("B149").Value = -9730
$ws.Range("C149").Value = 3187
$ws.Range("D149").Value = 854
$ws.Range("E149").Value = 792
$ws.Range("F149").Value = 1542

# Restore default (unstyled) formatting so the cells don't keep a custom style,
# matching the unstyled cells used throughout the rest of column A.
$ws.Range("A148:A149").ClearFormats()
